$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Value = "Luca Perenzoni"
$ws.Range("B64").Value = "Lorenzo Canali | CGB Gamberoni"
$ws.Range("C64").Value = "Andrea Roveda | Pinguini Trentini"
$ws.Range("D64").Value = "Luca Perenzoni | CGB Gamberoni"
$ws.Range("E64").Value = "Carlo Stedile | Mai una gioia"
$ws.Range("F64").Value = "Davide Bazzano | IMONTAGNA"
